# "still completed cpc for User1 QC"
#
# This script reproduces, via the Excel COM object model, the edits that were
# made to MainController.xlsx:
#   1. On MAIN_CONTROLLER: swap the RunStatus ("Y"/"N") flags for the
#      Calculator row (B2) and the CPC_1stTouchPoint_Approval row (B6) -
#      i.e. Calculator is no longer run ("N") and CPC_1stTouchPoint_Approval
#      is now run ("Y").
#   2. On DATASHEET: rename the test data files used for the FOS and the
#      CPC_1stTouchPoint_Approval rows (D3 / D6), and reduce the
#      ImplicityWait for the CPC_1stTouchPoint_Approval row (E6) from 20 to 10.
#   3. Update the active sheet/selection so MAIN_CONTROLLER becomes the
#      selected tab (it was DATASHEET before).

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("MAIN_CONTROLLER")
$wsData = $wb.Worksheets.Item("DATASHEET")

# --- MAIN_CONTROLLER: flip the run flags for Calculator (row 2) and
#     CPC_1stTouchPoint_Approval (row 6) ---
$wsMain.Range("B2").Value = "Y"
$wsMain.Range("B6").Value = "N"

# --- DATASHEET: update the test data file names and wait value ---
$wsData.Range("D3").Value = "FOS7_AfterPostSanction.xlsx"
$wsData.Range("D6").Value = "CPC_Dynamic.xlsx"
$wsData.Range("E6").Value = 10

# --- Selection / active tab bookkeeping ---
# Touch DATASHEET's selection first ...
$wsData.Range("G16").Select()

# ... then make MAIN_CONTROLLER the active sheet/selection last, so it ends
# up as the selected tab in the saved workbook.
$wsMain.Activate()
$wsMain.Range("H14").Select()
